$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# --- Rubric value edits (GEOMETRY section, rows 5/6/18/30/31/32/33/54) ---

# Row 5: mark Milestone I ("I") item as complete -> Student Confidence(X) = "X"
$ws.Range("F5").Value = "X"

# Row 6: assign this feature to Milestone II
$ws.Range("E6").Value = "II"

# Row 18: mark Milestone I item as complete -> "X"
$ws.Range("F18").Value = "X"

# Row 30: mark Milestone I item as complete -> "X"
$ws.Range("F30").Value = "X"

# Rows 31-33: reassign these features from Milestone I to Milestone II
$ws.Range("E31").Value = "II"
$ws.Range("E32").Value = "II"
$ws.Range("E33").Value = "II"

# Row 54: reassign from Milestone I (complete) to Milestone II (not yet complete)
$ws.Range("E54").Value = "II"
$ws.Range("F54").ClearContents()

# --- Restore the view/selection state ---
$ws.Range("E6").Select()
$excel.ActiveWindow.ScrollRow = 3
$excel.ActiveWindow.ScrollColumn = 1
